$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-29 23:02:58"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 23:02:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 23:02:58"
